$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Company name correction
$ws.Range("B3").Value = "کیمیا-ص. معدنی کیمیای زنجان گستران"

# Row 8: fiscal period headers (drop oldest quarter, shift left, append new quarter)
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Row 9: publish-date headers (drop oldest date, shift left, append new date)
$ws.Range("D9").Value = "1400-11-11 (4)"
$ws.Range("E9").Value = "1401-04-08 (9)"
$ws.Range("F9").Value = "1401-05-11 (4)"
$ws.Range("G9").Value = "1401-08-29 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-28 (8)"
$ws.Range("J9").Value = "1401-05-11 (2)"
$ws.Range("K9").Value = "1401-08-29 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-28"

# Row 11: فروش (Sales)
$ws.Range("D11").Value = 651854
$ws.Range("E11").Value = 2666846
$ws.Range("F11").Value = 1095405
$ws.Range("G11").Value = 2164223
$ws.Range("H11").Value = 1049140
$ws.Range("I11").Value = 2770179
$ws.Range("J11").Value = 1416204
$ws.Range("K11").Value = 2722839
$ws.Range("L11").Value = 1155263
$ws.Range("M11").Value = 3552373

# Row 12: بهای تمام شده کالای فروش رفته (COGS)
$ws.Range("D12").Value = -348388
$ws.Range("E12").Value = -1782479
$ws.Range("F12").Value = -711157
$ws.Range("G12").Value = -1405859
$ws.Range("H12").Value = -722924
$ws.Range("I12").Value = -2151446
$ws.Range("J12").Value = -982683
$ws.Range("K12").Value = -1798225
$ws.Range("L12").Value = -751341
$ws.Range("M12").Value = -2253260

# Row 13: سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 303466
$ws.Range("E13").Value = 884367
$ws.Range("F13").Value = 384248
$ws.Range("G13").Value = 758364
$ws.Range("H13").Value = 326216
$ws.Range("I13").Value = 618733
$ws.Range("J13").Value = 433521
$ws.Range("K13").Value = 924614
$ws.Range("L13").Value = 403922
$ws.Range("M13").Value = 1299113

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$ws.Range("D14").Value = -13622
$ws.Range("E14").Value = -38998
$ws.Range("F14").Value = -24333
$ws.Range("G14").Value = -110315
$ws.Range("H14").Value = -49285
$ws.Range("I14").Value = -133084
$ws.Range("J14").Value = -50408
$ws.Range("K14").Value = -110952
$ws.Range("L14").Value = -80705
$ws.Range("M14").Value = -120372

# Row 15: هزینه کاهش ارزش دریافتنی‌ها (Impairment expense)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense)
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 10500
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 10500
$ws.Range("H16").Value = 5250
$ws.Range("I16").Value = 5250
$ws.Range("J16").Value = 7350
$ws.Range("K16").Value = 7350
$ws.Range("L16").Value = 7350
$ws.Range("M16").Value = 900996

# Row 17: سود (زیان) عملیاتی (Operating profit)
$ws.Range("D17").Value = 289844
$ws.Range("E17").Value = 855869
$ws.Range("F17").Value = 359915
$ws.Range("G17").Value = 658549
$ws.Range("H17").Value = 282181
$ws.Range("I17").Value = 490899
$ws.Range("J17").Value = 390463
$ws.Range("K17").Value = 821012
$ws.Range("L17").Value = 330567
$ws.Range("M17").Value = 2079737

# Row 18: هزینه های مالی (Financial expenses)
$ws.Range("D18").Value = -6508
$ws.Range("E18").Value = -2414
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = -3646
$ws.Range("J18").Value = -6045
$ws.Range("K18").Value = -4511
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -2019

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense)
$ws.Range("D19").Value = 84055
$ws.Range("E19").Value = 605
$ws.Range("F19").Value = 129285
$ws.Range("G19").Value = 146947
$ws.Range("H19").Value = 6123
$ws.Range("I19").Value = -106
$ws.Range("J19").Value = 55880
$ws.Range("K19").Value = 103499
$ws.Range("L19").Value = 14579
$ws.Range("M19").Value = 12247

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit)
$ws.Range("D20").Value = 367391
$ws.Range("E20").Value = 854060
$ws.Range("F20").Value = 489200
$ws.Range("G20").Value = 805496
$ws.Range("H20").Value = 288304
$ws.Range("I20").Value = 487147
$ws.Range("J20").Value = 440298
$ws.Range("K20").Value = 920000
$ws.Range("L20").Value = 345146
$ws.Range("M20").Value = 2089965

# Row 21: مالیات (Tax)
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops)
$ws.Range("D22").Value = 367391
$ws.Range("E22").Value = 854060
$ws.Range("F22").Value = 489200
$ws.Range("G22").Value = 805496
$ws.Range("H22").Value = 288304
$ws.Range("I22").Value = 487147
$ws.Range("J22").Value = 440298
$ws.Range("K22").Value = 920000
$ws.Range("L22").Value = 345146
$ws.Range("M22").Value = 2089965

# Row 23: سود (زیان) عملیات متوقف شده (Discontinued ops)
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24: سود (زیان) خالص (Net profit)
$ws.Range("D24").Value = 367391
$ws.Range("E24").Value = 854060
$ws.Range("F24").Value = 489200
$ws.Range("G24").Value = 805496
$ws.Range("H24").Value = 288304
$ws.Range("I24").Value = 487147
$ws.Range("J24").Value = 440298
$ws.Range("K24").Value = 920000
$ws.Range("L24").Value = 345146
$ws.Range("M24").Value = 2089965

# Row 25: سود هر سهم پس از کسر مالیات (EPS after tax)
$ws.Range("D25").Value = 350
$ws.Range("E25").Value = 813
$ws.Range("F25").Value = 245
$ws.Range("G25").Value = 403
$ws.Range("H25").Value = 144
$ws.Range("I25").Value = 464
$ws.Range("J25").Value = 220
$ws.Range("K25").Value = 460
$ws.Range("L25").Value = 173
$ws.Range("M25").Value = 697

# Row 26: سرمایه (Capital)
$ws.Range("D26").Value = 1050000
$ws.Range("E26").Value = 1050000
$ws.Range("F26").Value = 2000000
$ws.Range("G26").Value = 2000000
$ws.Range("H26").Value = 2000000
$ws.Range("I26").Value = 1050000
$ws.Range("J26").Value = 2000000
$ws.Range("K26").Value = 2000000
$ws.Range("L26").Value = 2000000
$ws.Range("M26").Value = 3000000

# Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS based on latest capital)
$ws.Range("D27").Value = 122
$ws.Range("E27").Value = 285
$ws.Range("F27").Value = 163
$ws.Range("G27").Value = 269
$ws.Range("H27").Value = 96
$ws.Range("I27").Value = 162
$ws.Range("J27").Value = 147
$ws.Range("K27").Value = 307
$ws.Range("L27").Value = 115
$ws.Range("M27").Value = 697

